$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# li600 data upload: the weekly header row only had two dates (date_made,
# week 01). Extend it across the rest of the uploaded weeks -- C4 becomes
# "=B4+7" and D4:G4 are filled right with the same "+7 from the cell to
# the left" pattern, matching a drag-fill / Ctrl+R across the row.
$ws.Range("C4").Formula = "=B4+7"
$ws.Range("D4:G4").Formula = "=C4+7"

# New data now spans columns B:G, so auto-fit them to the (now uniform)
# content width instead of the old two custom column widths.
$ws.Range("B1:G9").Columns.AutoFit()
$ws.Range("B:G").ColumnWidth = 8.5

# Leave the cursor on the last cell of the newly-uploaded data.
$ws.Range("G9").Select() | Out-Null

$wb.Save() | Out-Null
